$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.21985324739876
$ws.Range("C2").Value = 0.1440013320705305
$ws.Range("D2").Value = 0.4993298101396988
$ws.Range("E2").Value = 0.1573178105596469
$ws.Range("G2").Value = 0.002537923202066159
$ws.Range("J2").Value = 0.05875730790749145
$ws.Range("L2").Value = 0.4867042824572394
$ws.Range("O2").Value = 6.669693685037032
$ws.Range("B3").Value = 2.101953421019232
$ws.Range("C3").Value = 0.1293848621976394
$ws.Range("D3").Value = 0.4988679542650516
$ws.Range("E3").Value = 0.1583032374668072
$ws.Range("G3").Value = 0.002541841759727814
$ws.Range("J3").Value = 0.05895870412865545
$ws.Range("L3").Value = 0.4780200569104238
$ws.Range("O3").Value = 6.707583413049122
$ws.Range("B4").Value = 2.030207835951217
$ws.Range("C4").Value = 0.1203647229137914
$ws.Range("D4").Value = 0.4988082704621633
$ws.Range("E4").Value = 0.1589670662510887
$ws.Range("G4").Value = 0.002544375977173839
$ws.Range("J4").Value = 0.05909110588833144
$ws.Range("L4").Value = 0.4728784955073735
$ws.Range("O4").Value = 6.734892475370373
$ws.Range("B5").Value = 2.001134670383124
$ws.Range("C5").Value = 0.1166776702344237
$ws.Range("D5").Value = 0.4988402887268535
$ws.Range("E5").Value = 0.1592523735416389
$ws.Range("G5").Value = 0.002545441031788535
$ws.Range("J5").Value = 0.05914726671269044
$ws.Range("L5").Value = 0.4708313199006398
$ws.Range("O5").Value = 6.747036837558539
$ws.Range("B6").Value = 1.996317026480142
$ws.Range("C6").Value = 0.1160647626438731
$ws.Range("D6").Value = 0.4988490092947302
$ws.Range("E6").Value = 0.1593006423584455
$ws.Range("G6").Value = 0.002545619839600633
$ws.Range("J6").Value = 0.05915672563132546
$ws.Range("L6").Value = 0.4704942940030747
$ws.Range("O6").Value = 6.749114705741022
$ws.Range("B7").Value = 2.029815080138974
$ws.Range("C7").Value = 0.1203150434253359
$ws.Range("D7").Value = 0.4988084741009402
$ws.Range("E7").Value = 0.1589708540971202
$ws.Range("G7").Value = 0.00254439020961384
$ws.Range("J7").Value = 0.05909185435089892
$ws.Range("L7").Value = 0.4728506918468298
$ws.Range("O7").Value = 6.735052147791436
$ws.Range("B8").Value = 2.179068263126283
$ws.Range("C8").Value = 0.1389711368525184
$ws.Range("D8").Value = 0.4991241141733127
$ws.Range("E8").Value = 0.1576453951242183
$ws.Range("G8").Value = 0.002539247772757327
$ws.Range("J8").Value = 0.05882493898685315
$ws.Range("L8").Value = 0.4836704770103353
$ws.Range("O8").Value = 6.681917779821049
$ws.Range("B9").Value = 2.476827536073927
$ws.Range("C9").Value = 0.1751878194824883
$ws.Range("D9").Value = 0.5015185338012174
$ws.Range("E9").Value = 0.1555119702967751
$ws.Range("G9").Value = 0.002530176077038834
$ws.Range("J9").Value = 0.05837057019032166
$ws.Range("L9").Value = 0.5063963179071607
$ws.Range("O9").Value = 6.609879168853723
$ws.Range("B10").Value = 2.698647636463136
$ws.Range("C10").Value = 0.2015662281255004
$ws.Range("D10").Value = 0.5043595780851007
$ws.Range("E10").Value = 0.1542278167996525
$ws.Range("G10").Value = 0.002524121884755595
$ws.Range("J10").Value = 0.05807840011532051
$ws.Range("L10").Value = 0.5240095909055782
$ws.Range("O10").Value = 6.576653749021318
$ws.Range("B11").Value = 2.800217185582198
$ws.Range("C11").Value = 0.213515519869901
$ws.Range("D11").Value = 0.5058869331873836
$ws.Range("E11").Value = 0.1537049924715639
$ws.Range("G11").Value = 0.00252149890999145
$ws.Range("J11").Value = 0.0579544369440459
$ws.Range("L11").Value = 0.5322209379862244
$ws.Range("O11").Value = 6.565838098243603
$ws.Range("B12").Value = 2.838773204413314
$ws.Range("C12").Value = 0.2180330254574017
$ws.Range("D12").Value = 0.5064990703213681
$ws.Range("E12").Value = 0.1535158225729134
$ws.Range("G12").Value = 0.002520524405396758
$ws.Range("J12").Value = 0.05790877458350963
$ws.Range("L12").Value = 0.5353588936927878
$ws.Range("O12").Value = 6.562362321563796
$ws.Range("B13").Value = 2.830465325002194
$ws.Range("C13").Value = 0.2170604332789026
$ws.Range("D13").Value = 0.5063657344946648
$ws.Range("E13").Value = 0.153556171903217
$ws.Range("G13").Value = 0.002520733449642169
$ws.Range("J13").Value = 0.05791855197592
$ws.Range("L13").Value = 0.5346818132124156
$ws.Range("O13").Value = 6.563083295670026
$ws.Range("B14").Value = 2.803387339942731
$ws.Range("C14").Value = 0.2138873277785649
$ws.Range("D14").Value = 0.5059366175969728
$ws.Range("E14").Value = 0.153689252787899
$ws.Range("G14").Value = 0.002521418361561767
$ws.Range("J14").Value = 0.05795065466293714
$ws.Range("L14").Value = 0.5324785290307545
$ws.Range("O14").Value = 6.565539710904147
$ws.Range("B15").Value = 2.786813470414529
$ws.Range("C15").Value = 0.2119427356921335
$ws.Range("D15").Value = 0.5056781670580364
$ws.Range("E15").Value = 0.1537719160387052
$ws.Range("G15").Value = 0.002521840330051099
$ws.Range("J15").Value = 0.05797048495659318
$ws.Range("L15").Value = 0.5311326615604202
$ws.Range("O15").Value = 6.567125113395548
$ws.Range("B16").Value = 2.692023035503269
$ws.Range("C16").Value = 0.2007842839474279
$ws.Range("D16").Value = 0.5042644875068873
$ws.Range("E16").Value = 0.1542632182326464
$ws.Range("G16").Value = 0.002524295932275802
$ws.Range("J16").Value = 0.05808668085229662
$ws.Range("L16").Value = 0.5234769540563775
$ws.Range("O16").Value = 6.577447226061992
$ws.Range("B17").Value = 2.634040824240287
$ws.Range("C17").Value = 0.1939258974310007
$ws.Range("D17").Value = 0.5034573994111042
$ws.Range("E17").Value = 0.1545803214631931
$ws.Range("G17").Value = 0.00252583587733545
$ws.Range("J17").Value = 0.0581602498891538
$ws.Range("L17").Value = 0.5188313056064544
$ws.Range("O17").Value = 6.584881744329323
$ws.Range("B18").Value = 2.600753428533949
$ws.Range("C18").Value = 0.1899764039935405
$ws.Range("D18").Value = 0.5030152996945674
$ws.Range("E18").Value = 0.1547684853158326
$ws.Range("G18").Value = 0.002526733959098246
$ws.Range("J18").Value = 0.05820340733109752
$ws.Range("L18").Value = 0.5161779873179739
$ws.Range("O18").Value = 6.589562398869305
$ws.Range("B19").Value = 2.589493656174113
$ws.Range("C19").Value = 0.1886383666140148
$ws.Range("D19").Value = 0.502869411534931
$ws.Range("E19").Value = 0.154833186403037
$ws.Range("G19").Value = 0.002527040157327502
$ws.Range("J19").Value = 0.05821816460110618
$ws.Range("L19").Value = 0.5152828408191823
$ws.Range("O19").Value = 6.591216611904912
$ws.Range("B20").Value = 2.640206677879348
$ws.Range("C20").Value = 0.1946564753869495
$ws.Range("D20").Value = 0.5035410265739699
$ws.Range("E20").Value = 0.1545459677119148
$ws.Range("G20").Value = 0.002525670670338931
$ws.Range("J20").Value = 0.05815233119500895
$ws.Range("L20").Value = 0.5193239045796219
$ws.Range("O20").Value = 6.584048448787513
$ws.Range("B21").Value = 2.811338265121492
$ws.Range("C21").Value = 0.2148195489551767
$ws.Range("D21").Value = 0.5060617435631229
$ws.Range("E21").Value = 0.1536499246115444
$ws.Range("G21").Value = 0.00252121667785404
$ws.Range("J21").Value = 0.05794119064513481
$ws.Range("L21").Value = 0.5331249145246062
$ws.Range("O21").Value = 6.564801365222877
$ws.Range("B22").Value = 2.923728949298891
$ws.Range("C22").Value = 0.2279538802532954
$ws.Range("D22").Value = 0.5079059514056894
$ws.Range("E22").Value = 0.1531156696536957
$ws.Range("G22").Value = 0.002518415033861803
$ws.Range("J22").Value = 0.05781065471095026
$ws.Range("L22").Value = 0.5423107211064035
$ws.Range("O22").Value = 6.555836195593884
$ws.Range("B23").Value = 2.863694404024614
$ws.Range("C23").Value = 0.2209478727637588
$ws.Range("D23").Value = 0.5069036660817972
$ws.Range("E23").Value = 0.1533961149143952
$ws.Range("G23").Value = 0.002519900354462147
$ws.Range("J23").Value = 0.05787964410304802
$ws.Range("L23").Value = 0.5373929275986455
$ws.Range("O23").Value = 6.560289832270257
$ws.Range("B24").Value = 2.637418947849028
$ws.Range("C24").Value = 0.1943262013772369
$ws.Range("D24").Value = 0.5035031504980054
$ws.Range("E24").Value = 0.1545614807981455
$ws.Range("G24").Value = 0.002525745320674586
$ws.Range("J24").Value = 0.05815590855188013
$ws.Range("L24").Value = 0.5191011459837114
$ws.Range("O24").Value = 6.58442391570739
$ws.Range("B25").Value = 2.395736467477889
$ws.Range("C25").Value = 0.1654302363887723
$ws.Range("D25").Value = 0.5006807602118926
$ws.Range("E25").Value = 0.1560393231853539
$ws.Range("G25").Value = 0.002532522481718646
$ws.Range("J25").Value = 0.05848614244936101
$ws.Range("L25").Value = 0.5346818132124156
$ws.Range("O25").Value = 6.625914799076526
